$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the column header: "question_id" -> "question-id" (fixing the bedrock model)
$ws.Range("A1").Value = "question-id"

# Match the author's cursor/selection position (B8) as seen in the diff
$ws.Range("B8").Select()
